$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "removed redundant batAve calcs in setup"
#
# Two new benchmark checkpoint sheets are appended (each a copy of the sheet
# before it, carrying the same layout/formulas but fresh timing numbers):
#   16. "optimized rId in participants"  (copy of "created c_did_get_hit")
#   17. "no redun batAve calcs"          (copy of "optimized rId in participants")
# ---------------------------------------------------------------------------

$prev = $wb.Worksheets.Item("created c_did_get_hit")

# --- Sheet 16: "optimized rId in participants" ------------------------------
$prev.Copy([System.Reflection.Missing]::Value, $prev)
$s16 = $wb.Worksheets.Item("created c_did_get_hit (2)")
$s16.Name = "optimized rId in participants"

$s16.Range("B2").Value = 5.612
$s16.Range("C2").Value = 5.733
$s16.Range("D2").Value = 5.663
$s16.Range("F2").Formula = "=('created c_did_get_hit'!E2 - E2)/'created c_did_get_hit'!E2"

$s16.Range("B3").Value = 4.924
$s16.Range("C3").Value = 5.036
$s16.Range("D3").Value = 4.976
$s16.Range("H3").Value = "We're now returning a set instead of a generator from get_participants, "

$s16.Range("B4").Value = 0.667
$s16.Range("C4").Value = 0.696
$s16.Range("D4").Value = 0.686
$s16.Range("H4").Value = "which allows for hash-based lookups"

$s16.Activate()
$s16.Range("A1:H6").Select()

# --- Sheet 17: "no redun batAve calcs" --------------------------------------
$s16.Copy([System.Reflection.Missing]::Value, $s16)
$s17 = $wb.Worksheets.Item("optimized rId in participants (2)")
$s17.Name = "no redun batAve calcs"

$s17.Range("B2").Value = 3.633
$s17.Range("C2").Value = 3.563
$s17.Range("D2").Value = 3.696
$s17.Range("F2").Formula = "=('optimized rId in participants'!E2-E2)/'optimized rId in participants'!E2"

$s17.Range("B3").Value = 3.241
$s17.Range("C3").Value = 3.185
$s17.Range("D3").Value = 3.289
$s17.Range("H3").Value = "Setup was having the players calculate their batting averages again…"

$s17.Range("B4").Value = 0.291
$s17.Range("C4").Value = 0.378
$s17.Range("D4").Value = 0.394
$s17.Range("H4").Value = "even though the batave csvs have that info!! "

$s17.Range("H5").Value = "Dumb."

$s17.Activate()
$s17.Range("H5").Select()

# --- Tidy up the no-longer-active previous sheet's selection ---------------
$prev.Range("A1:H6").Select()
$s17.Activate()
